# Apply updated crypto price / 1h-volume-change figures (refreshed symbol list run).
# Source values taken verbatim from the diff; cells are plain text (inlineStr) in the
# original workbook, so we force a Text number format before assigning the value and
# then clear the format back off, which keeps Excel from "helpfully" re-interpreting
# strings like "302.88", "0.07896" or "1.15%" as numbers/percentages (which would drop
# significant trailing zeros or convert "%" strings into fractional numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '302.88'
    'E2' = '1.15%'
    'D3' = '31.87'
    'E3' = '0.68%'
    'D4' = '4.997'
    'E4' = '-1.95%'
    'D5' = '0.07896'
    'E5' = '-3.61%'
    'D6' = '2.116'
    'E6' = '-18.39%'
    'D7' = '7.853'
    'E7' = '0.14%'
    'E8' = '-1.20%'
    'D9' = '0.9294'
    'E9' = '0.49%'
    'D10' = '0.1754'
    'E10' = '-0.12%'
    'D11' = '0.07968'
    'E11' = '6.75%'
    'D12' = '0.08760'
    'E12' = '-2.25%'
    'D13' = '0.03133'
    'E13' = '4.72%'
    'D14' = '0.1002'
    'E14' = '0.43%'
    'D15' = '0.001509'
    'E15' = '-0.27%'
    'D16' = '0.005957'
    'E16' = '1.47%'
    'D17' = '3.468'
    'E17' = '-3.45%'
    'E18' = '0.86%'
    'D19' = '0.3294'
    'E19' = '1.51%'
    'E20' = '-4.09%'
    'D21' = '4.152'
    'E21' = '6.33%'
    'D22' = '0.1791'
    'E22' = '5.57%'
    'D23' = '0.04603'
    'E23' = '-0.03%'
    'D24' = '0.001236'
    'E24' = '-0.64%'
    'D25' = '0.004504'
    'E25' = '-0.62%'
    'D26' = '0.0001248'
    'E26' = '4.26%'
    'D39' = '0.01733'
    'E39' = '-2.28%'
    'D40' = '0.04767'
    'E40' = '5.20%'
    'D41' = '0.007386'
    'E41' = '5.71%'
    'E42' = '-0.88%'
    'D43' = '0.002347'
    'E43' = '9.85%'
    'D44' = '0.01101'
    'E44' = '14.58%'
    'D45' = '0.00006042'
    'E45' = '-3.93%'
    'D46' = '0.00000000749'
    'E46' = '0.13%'
    'D47' = '0.003397'
    'E47' = '-61.13%'
    'D48' = '0.8204'
    'E48' = '0.60%'
    'D49' = '0.00002098'
    'E49' = '0.13%'
    'D50' = '0.0001998'
    'E50' = '0.13%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
